$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.174232602119446
$ws.Range("B1").Value = 2.17740797996521
$ws.Range("C1").Value = 10.43294429779053
$ws.Range("D1").Value = 2.562594652175903
$ws.Range("E1").Value = 1.249155879020691
